$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.916.49"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").Value = "1.554.51"
$ws.Range("E3").Value = "  +0.74%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "206.81"
$ws.Range("E5").Value = "  +0.37%  "
$ws.Range("D6").Value = "0.489"
$ws.Range("E6").Value = "  +0.36%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").Value = "21.92"
$ws.Range("E8").Value = "  +2.35%  "
$ws.Range("E9").Value = "  +0.25%  "
$ws.Range("D10").Value = "0.0586"
$ws.Range("E10").Value = "  +0.78%  "
$ws.Range("E11").Value = "  +0.62%  "
$ws.Range("D12").Value = "1.777.03"
$ws.Range("E12").Value = "  +0.82%  "
$ws.Range("D13").Value = "1.555.64"
$ws.Range("E13").Value = "  +0.81%  "
$ws.Range("E14").Value = "  +1.58%  "
$ws.Range("D15").Value = "0.517"
$ws.Range("E15").Value = "  +1.46%  "
$ws.Range("D16").Value = "26.926.83"
$ws.Range("E16").Value = "  +0.38%  "
$ws.Range("D17").Value = "61.69"
$ws.Range("E17").Value = "  +0.70%  "
$ws.Range("D18").Value = "217.67"
$ws.Range("E18").Value = "  +1.20%  "
$ws.Range("D19").Value = "0.0₃0693"
$ws.Range("E19").Value = "  +1.59%  "
$ws.Range("D20").Value = "7.29"
$ws.Range("E20").Value = "  +0.77%  "
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("E22").Value = "  +1.15%  "
$ws.Range("D23").Value = "9.19"
$ws.Range("E23").Value = "  +0.44%  "
$ws.Range("E24").Value = "  +0.65%  "
$ws.Range("D25").Value = "154.13"
$ws.Range("E25").Value = "  +0.89%  "
$ws.Range("D26").Value = "6.59"
$ws.Range("E26").Value = "  -0.29%  "
$ws.Range("D27").Value = "14.91"
$ws.Range("E27").Value = "  +0.45%  "
$ws.Range("E28").Value = "  +0.13%  "
$ws.Range("E29").Value = "  +0.72%  "
$ws.Range("D30").Value = "0.0468"
$ws.Range("E30").Value = "  +2.11%  "
$ws.Range("E31").Value = "  +0.18%  "
$ws.Range("E32").Value = "  -0.07%  "
$ws.Range("D33").Value = "1.436.90"
$ws.Range("E33").Value = "  +5.08%  "
$ws.Range("D34").Value = "3.06"
$ws.Range("E34").Value = "  +3.89%  "
$ws.Range("E35").Value = "  +3.51%  "
$ws.Range("D36").Value = "0.972"
$ws.Range("E36").Value = "  +0.84%  "
$ws.Range("E37").Value = "  +0.73%  "
$ws.Range("D38").Value = "0.0164"
$ws.Range("E38").Value = "  +0.05%  "
$ws.Range("D39").Value = "0.519"
$ws.Range("E39").Value = "  -0.27%  "
$ws.Range("D40").Value = "0.811"
$ws.Range("E40").Value = "  +0.50%  "
$ws.Range("E41").Value = "  +0.16%  "
$ws.Range("E42").Value = "  -1.57%  "
$ws.Range("D43").Value = "0.987"
$ws.Range("E43").Value = "  -0.35%  "
$ws.Range("D44").Value = "2.27"
$ws.Range("E44").Value = "  +2.66%  "
$ws.Range("D45").Value = "63.97"
$ws.Range("E45").Value = "  +1.23%  "
$ws.Range("D46").Value = "1.76"
$ws.Range("E46").Value = "  +1.50%  "
$ws.Range("D47").Value = "1.690.98"
$ws.Range("E47").Value = "  +0.82%  "
$ws.Range("E48").Value = "  +2.68%  "
$ws.Range("E49").Value = "  +2.72%  "
$ws.Range("D50").Value = "0.0₆0100"
$ws.Range("E50").Value = "  +3.71%  "
$ws.Range("D51").Value = "0.0955"
$ws.Range("E51").Value = "  +1.71%  "
